# Generate Report for Handoff
# Update the "Latest Handoff" timestamps for the file
# 4faca4ca-afda-4085-8a9d-7691c44b8dc4.md now that it has been handed off again.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest Handoff Date" column (D) for the 4faca4ca row (row 6)
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-03-22 06:41:37"

# zh-cn sheet: "Latest Handoff Datetime" column (E) for the 4faca4ca row (row 6)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-22 06:41:30"

# de-de sheet: "Latest Handoff Datetime" column (E) for the 4faca4ca row (row 6)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-22 06:41:37"
